$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- SARIF2001 (AuthorHighQualityMessages) notes block, rows 37-39 ---
# These cells move from "IN PROGRESS: LJG" (style s=12, gold fill) to
# "READY FOR REVIEW" (style s=7, blue fill) -- copy both value + format
# from an existing "READY FOR REVIEW" cell (G2) so the style index matches
# exactly what the rest of the sheet uses for that status.
$ws.Range("G2").Copy($ws.Range("G37"))
$ws.Range("G2").Copy($ws.Range("G38"))
$ws.Range("G2").Copy($ws.Range("G39"))

# New TODO notes for the two new checks (H38/H39), styled like the existing
# note in H62 (s=16).
$ws.Range("H62").Copy($ws.Range("H38"))
$ws.Range("H38").Value = "TODO: IncludeDynamicContent applies to Markdown. Enquote might not because you can set values of with code font."

$ws.Range("H62").Copy($ws.Range("H39"))
$ws.Range("H39").Value = "TODO: Don't include the string itself. DO include the rule id."

# --- SARIF2004 (OptimizeFileSize) rows 49-50 ---
# F49 goes from "IN PROGRESS: HK" (s=12) to "DONE" (s=11, green fill) --
# copy format from another DONE cell (F2).
$ws.Range("F2").Copy($ws.Range("F49"))
$ws.Range("F49").Value = "DONE"

# F50 stays the same gold "IN PROGRESS"-style cell, just the text changes
# from "IN PROGRESS: HK" to "IN PROGRESS".
$ws.Range("F50").Value = "IN PROGRESS"

# --- SARIF2005 rule rename ---
# B52: "ProvideHelpfulToolInformation" -> "ProvideToolProperties"
$ws.Range("B52").Value = "ProvideToolProperties"

# --- sheet view / selection state ---
# (topLeftCell scroll position isn't modeled by this engine's object model,
# but the active selection is.)
$ws.Range("I3").Select()
